$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.740.34"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "3.812.68"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'710.07"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'170.08"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("D7").Value = "3.812.91"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("D11").Value = "'7.54"
$ws.Range("E11").Value = "  +4.78%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "'35.94"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").Value = "4.457.45"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").Value = "3.802.16"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "70.825.69"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'7.14"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").Value = "'17.27"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").Value = "'495.77"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'10.62"
$ws.Range("E22").Value = "  -4.86%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'84.19"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").Value = "'12.09"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("D28").Value = "3.964.62"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -5.20%  "
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("D34").Value = "'29.10"
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("D36").Value = "3.782.34"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").Value = "'9.10"
$ws.Range("E37").Value = "  -2.34%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'0.101"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("D42").Value = "'5.96"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("D43").Value = "'3.23"
$ws.Range("E43").Value = "  -4.91%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "'0.000323"
$ws.Range("E46").Value = "  +4.46%  "
$ws.Range("D47").Value = "'165.48"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "'48.82"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "'423.01"
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("D50").Value = "'8.62"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("E51").Value = "  -3.58%  "
